$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.838.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.422.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.418.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  -8.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.408"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.011.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.929.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("E17").Value = "  -4.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.421.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -5.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.55%  "
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.523"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.05%  "
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("E28").Value = "  -6.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.177"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  -8.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.80%  "
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.685.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0679"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0287"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "304.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.810"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.94%  "
